$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 2) under the existing header row.
$ws.Cells.Item(2, 1).Value = "HOSP..DE.REHABILITACION.PSICOFISICA"
$ws.Cells.Item(2, 2).Value = "HOSP..DE.REHABILITACION.PSICOFISICA"
# Column C (ID) has no value for this record - touch formatting only so the
# cell exists in the sheet without holding any content.
$ws.Cells.Item(2, 3).Borders.LineStyle = -4142
$ws.Cells.Item(2, 4).Value = "IREP"
